$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Radio-button alternatives example: swap the symbol set ---
# Row 19 (A19): was the "◎" glyph at 16pt Arial; now shows the filled "◉"
# glyph at a smaller 12pt Calibri (theme minor) weight.
$ws.Range("A19").Value = "◉"
$ws.Range("A19").Font.Name = "Calibri"
$ws.Range("A19").Font.ThemeFont = 1   # xlThemeFontMinor
$ws.Range("A19").Font.Size = 12
$ws.Rows.Item(19).AutoFit()

# Row 20 (A20): was the "◉" glyph at 22pt Arial; now shows the open-circle
# "○" glyph, switched to Calibri (theme minor) but keeping the 22pt size.
$ws.Range("A20").Value = "○"
$ws.Range("A20").Font.Name = "Calibri"
$ws.Range("A20").Font.ThemeFont = 1   # xlThemeFontMinor
$ws.Rows.Item(20).RowHeight = 29

# Row 21 (A21) and Row 22 (A22): the extra "🔘" / "⊙" alternatives are
# dropped from this example - clear their text but keep the existing look.
$ws.Range("A21").Value = ""
$ws.Range("A22").Value = ""

# Selection cursor moved down to where the edits were made.
$ws.Range("C18").Select()
